$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Daily refresh of the cryptocurrency price/volume table (GitHub Actions job).
# Columns D (Price) and E (Volume 1h) are stored as text in this sheet, and some
# of the new values look like numbers (e.g. "1.001"), so COM would otherwise
# auto-convert them. Force the whole data range to text format before writing,
# then restore the default "Normal" style afterwards so no stray number format
# is left behind on the cells.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '20.196.25'
$ws.Range("E2").Value = '  +1.05%  '
$ws.Range("D3").Value = '1.438.47'
$ws.Range("E3").Value = '  +1.73%  '
$ws.Range("D4").Value = '1.009'
$ws.Range("E4").Value = '  +1.03%  '
$ws.Range("D5").Value = '0.9217'
$ws.Range("E5").Value = '  -7.76%  '
$ws.Range("D6").Value = '274.32'
$ws.Range("E6").Value = '  -0.02%  '
$ws.Range("D7").Value = '0.3638'
$ws.Range("E7").Value = '  -1.41%  '
$ws.Range("D8").Value = '0.3051'
$ws.Range("E8").Value = '  -2.21%  '
$ws.Range("D9").Value = '38.67'
$ws.Range("E9").Value = '  -2.32%  '
$ws.Range("D10").Value = '1.011'
$ws.Range("E10").Value = '  -2.54%  '
$ws.Range("D11").Value = '0.06439'
$ws.Range("D12").Value = '1.001'
$ws.Range("E12").Value = '  +0.14%  '
$ws.Range("D13").Value = '5.288'
$ws.Range("E13").Value = '  -3.12%  '
$ws.Range("D14").Value = '17.22'
$ws.Range("E14").Value = '  -2.09%  '
$ws.Range("D15").Value = '5.995'
$ws.Range("E15").Value = '  -3.01%  '
$ws.Range("D16").Value = '0.000009987'
$ws.Range("E16").Value = '  -1.59%  '
$ws.Range("D17").Value = '1.441.10'
$ws.Range("E17").Value = '  +1.88%  '
$ws.Range("D18").Value = '0.9399'
$ws.Range("E18").Value = '  -5.92%  '
$ws.Range("D19").Value = '0.05620'
$ws.Range("E19").Value = '  -1.15%  '
$ws.Range("D20").Value = '67.54'
$ws.Range("E20").Value = '  -4.50%  '
$ws.Range("D21").Value = '5.307'
$ws.Range("E21").Value = '  -5.01%  '
$ws.Range("D22").Value = '14.09'
$ws.Range("E22").Value = '  -4.35%  '
$ws.Range("D23").Value = '10.68'
$ws.Range("E23").Value = '  -3.58%  '
$ws.Range("D24").Value = '2.241'
$ws.Range("E24").Value = '  -1.56%  '
$ws.Range("D25").Value = '20.249.91'
$ws.Range("E25").Value = '  +1.24%  '
$ws.Range("D26").Value = '138.84'
$ws.Range("E26").Value = '  +2.55%  '
$ws.Range("D27").Value = '2.022'
$ws.Range("E27").Value = '  -10.59%  '
$ws.Range("D28").Value = '16.83'
$ws.Range("E28").Value = '  -0.97%  '
$ws.Range("D29").Value = '1.592.21'
$ws.Range("E29").Value = '  +1.33%  '
$ws.Range("D30").Value = '110.04'
$ws.Range("E30").Value = '  +0.63%  '
$ws.Range("D31").Value = '4.035'
$ws.Range("E31").Value = '  -1.29%  '
$ws.Range("D32").Value = '4.759'
$ws.Range("E32").Value = '  -10.74%  '
$ws.Range("D33").Value = '0.07636'
$ws.Range("E33").Value = '  -0.72%  '
$ws.Range("D34").Value = '0.7694'
$ws.Range("E34").Value = '  -6.30%  '
$ws.Range("D35").Value = '1.456'
$ws.Range("E35").Value = '  +0.33%  '
$ws.Range("D36").Value = '0.05674'
$ws.Range("E36").Value = '  -2.95%  '
$ws.Range("D37").Value = '1.126'
$ws.Range("E37").Value = '  +3.29%  '
$ws.Range("D38").Value = '4.597'
$ws.Range("E38").Value = '  -6.11%  '
$ws.Range("B39").Value = 'Frax'
$ws.Range("C39").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D39").Value = '0.9355'
$ws.Range("E39").Value = '  -6.35%  '
$ws.Range("B40").Value = 'VeChain'
$ws.Range("C40").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D40").Value = '0.01965'
$ws.Range("E40").Value = '  -4.72%  '
$ws.Range("E41").Value = '  -3.62%  '
$ws.Range("D42").Value = '0.1820'
$ws.Range("E42").Value = '  -4.70%  '
$ws.Range("D43").Value = '6.891'
$ws.Range("E43").Value = '  -17.96%  '
$ws.Range("D44").Value = '3.470'
$ws.Range("E44").Value = '  -1.56%  '
$ws.Range("D45").Value = '0.5144'
$ws.Range("E45").Value = '  -2.78%  '
$ws.Range("D46").Value = '11.62'
$ws.Range("D47").Value = '114.26'
$ws.Range("E47").Value = '  +0.78%  '
$ws.Range("D48").Value = '0.5019'
$ws.Range("E48").Value = '  -2.58%  '
$ws.Range("D49").Value = '1.713'
$ws.Range("E49").Value = '  -3.04%  '
$ws.Range("E50").Value = '  +2.67%  '
$ws.Range("D51").Value = '0.9846'
$ws.Range("E51").Value = '  -1.51%  '

$ws.Range("D2:E51").Style = "Normal"
